$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (want-to-go count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 23
$wsExhibit.Range("F9").Value = 314

# Sheet "全部类型" (All types) mirrors the same data - apply identical updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 23
$wsAll.Range("F9").Value = 314
